$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("C1").Value = "First Review(25%)"
$ws.Range("D1").Value = "Second Review (35%)"
$ws.Columns.Item(3).ColumnWidth = 24
